$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (A1:D1) to new snake_case column codes
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'

# Title-case the connector words (de/del/la/las/el/los/y) in state/municipality names,
# plus a couple of one-off casing fixes (GUANAJUATO -> Guanajuato, MonteMorelos -> Montemorelos)
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B7').Value = 'San Francisco De Los Romo'
$ws.Range('B11').Value = 'Playas De Rosarito'
$ws.Range('B24').Value = 'Amatenango De La Frontera'
$ws.Range('B32').Value = 'Comitán De Domínguez'
$ws.Range('B53').Value = 'Ocozocoautla De Espinosa'
$ws.Range('B58').Value = 'San Cristóbal De Las Casas'
$ws.Range('B96').Value = 'Guadalupe Y Calvo'
$ws.Range('B99').Value = 'Hidalgo Del Parral'
$ws.Range('B114').Value = 'San Francisco Del Oro'
$ws.Range('B121').Value = 'Valle De Zaragoza'
$ws.Range('B140').Value = 'Villa De Álvarez'
$ws.Range('A142').Value = 'Ciudad De México'
$ws.Range('B146').Value = 'Cuajimalpa De Morelos'
$ws.Range('B159').Value = 'Coneto De Comonfort'
$ws.Range('B172').Value = 'Nombre De Dios'
$ws.Range('B176').Value = 'Pánuco De Coronado'
$ws.Range('B183').Value = 'San Pedro Del Gallo'
$ws.Range('A193').Value = 'Estado De México'
$ws.Range('B193').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B195').Value = 'Almoloya De Alquisiras'
$ws.Range('B196').Value = 'Almoloya De Juárez'
$ws.Range('B200').Value = 'Atizapán De Zaragoza'
$ws.Range('B205').Value = 'Chapa De Mota'
$ws.Range('B209').Value = 'Coacalco De Berriozábal'
$ws.Range('B212').Value = 'Ecatepec De Morelos'
$ws.Range('B215').Value = 'Ixtapan De La Sal'
$ws.Range('B216').Value = 'Ixtapan Del Oro'
$ws.Range('B224').Value = 'Naucalpan De Juárez'
$ws.Range('B230').Value = 'San Felipe Del Progreso'
$ws.Range('B238').Value = 'Tenango Del Valle'
$ws.Range('B245').Value = 'Tlalnepantla De Baz'
$ws.Range('B249').Value = 'Valle De Bravo'
$ws.Range('B250').Value = 'Villa De Allende'
$ws.Range('B261').Value = 'Apaseo El Alto'
$ws.Range('B262').Value = 'Apaseo El Grande'
$ws.Range('B269').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B273').Value = 'Jaral Del Progreso'
$ws.Range('B281').Value = 'Purísima Del Rincón'
$ws.Range('B286').Value = 'San Francisco Del Rincón'
$ws.Range('B288').Value = 'San Luis De La Paz'
$ws.Range('B289').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B290').Value = 'Silao De La Victoria'
$ws.Range('B294').Value = 'Valle De Santiago'
$ws.Range('B298').Value = 'Acapulco De Juárez'
$ws.Range('B300').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B301').Value = 'Alcozauca De Guerero'
$ws.Range('B305').Value = 'Atoyac De Álvarez'
$ws.Range('B306').Value = 'Ayutla De Los Libres'
$ws.Range('B309').Value = 'Chilapa De Álvarez'
$ws.Range('B310').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B311').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B316').Value = 'Coyuca De Benítez'
$ws.Range('B317').Value = 'Coyuca De Catalán'
$ws.Range('B320').Value = 'Cuetzala Del Progreso'
$ws.Range('B324').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B325').Value = 'Iguala De La Independencia'
$ws.Range('B327').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B328').Value = 'Zihuatanejo De Azueta'
$ws.Range('B330').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B333').Value = 'Mártir De Cuilapan'
$ws.Range('B346').Value = 'Taxco De Alarcón'
$ws.Range('B348').Value = 'Técpan De Galeana'
$ws.Range('B350').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B352').Value = 'Tixtla De Guerero'
$ws.Range('B355').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B369').Value = 'Cuautepec De Hinojosa'
$ws.Range('B372').Value = 'Huasca De Ocampo'
$ws.Range('B380').Value = 'Mineral Del Monte'
$ws.Range('B381').Value = 'Mixquiahuala De Juárez'
$ws.Range('B383').Value = 'Omitlán De Juárez'
$ws.Range('B384').Value = 'Pachuca De Soto'
$ws.Range('B388').Value = 'Santiago De Anaya'
$ws.Range('B391').Value = 'Tenango De Doria'
$ws.Range('B392').Value = 'Tepehuacán De Guerero'
$ws.Range('B393').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B395').Value = 'Tezontepec De Aldama'
$ws.Range('B401').Value = 'Tula De Allende'
$ws.Range('B402').Value = 'Tulancingo De Bravo'
$ws.Range('B405').Value = 'Zacualtipán De Ángeles'
$ws.Range('B406').Value = 'Zapotlán De Juárez'
$ws.Range('B408').Value = 'Ahualulco De Mercado'
$ws.Range('B412').Value = 'Atotonilco El Alto'
$ws.Range('B413').Value = 'Autlán De Navarro'
$ws.Range('B421').Value = 'Cuautitlán De García Barragán'
$ws.Range('B431').Value = 'Ixtlahuacán Del Río'
$ws.Range('B437').Value = 'Lagos De Moreno'
$ws.Range('B443').Value = 'Ojuelos De Jalisco'
$ws.Range('B448').Value = 'San Cristóbal De La Barranca'
$ws.Range('B449').Value = 'San Diego De Alejandría'
$ws.Range('B451').Value = 'San Juan De Los Lagos'
$ws.Range('B454').Value = 'San Martín De Bolaños'
$ws.Range('B455').Value = 'San Miguel El Alto'
$ws.Range('B456').Value = 'San Sebastián Del Oeste'
$ws.Range('B458').Value = 'Tamazula De Gordiano'
$ws.Range('B463').Value = 'Tepatitlán De Morelos'
$ws.Range('B465').Value = 'Tizapán El Alto'
$ws.Range('B466').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B475').Value = 'Unión De San Antonio'
$ws.Range('B476').Value = 'Unión De Tula'
$ws.Range('B480').Value = 'Yahualica De González Gallo'
$ws.Range('B481').Value = 'Zacoalco De Torres'
$ws.Range('B484').Value = 'Zapotlán El Grande'
$ws.Range('B499').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B553').Value = 'Tiquicheo De Nicolás Romero'
$ws.Range('B578').Value = 'Puente De Ixtla'
$ws.Range('B584').Value = 'Tlaltizapán De Zapata'
$ws.Range('B593').Value = 'Amatlán De Cañas'
$ws.Range('B596').Value = 'Ixtlán Del Río'
$ws.Range('B620').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B624').Value = 'Cuilápam De Guerero'
$ws.Range('B625').Value = 'El Barrio De La Soledad'
$ws.Range('B626').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B627').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B628').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B629').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B633').Value = 'Mazatlán Villa De Flores'
$ws.Range('B634').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B635').Value = 'Oaxaca De Juárez'
$ws.Range('B636').Value = 'Ocotlán De Morelos'
$ws.Range('B637').Value = 'Putla Villa De Guerero'
$ws.Range('B642').Value = 'San Baltazar Yatzachi El Bajo'
$ws.Range('B671').Value = 'San Pedro El Alto'
$ws.Range('B676').Value = 'San Pedro Y San Pablo Ayutla'
$ws.Range('B697').Value = 'Tlacolula De Matamoros'
$ws.Range('B698').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B699').Value = 'Villa De Zaachila'
$ws.Range('B700').Value = 'Villa Sola De Vega'
$ws.Range('B701').Value = 'Zimatlán De Álvarez'
$ws.Range('B712').Value = 'Chalchicomula De Sesma'
$ws.Range('B723').Value = 'Huehuetlán El Grande'
$ws.Range('B732').Value = 'Mazapiltepec De Juárez'
$ws.Range('B737').Value = 'Palmar De Bravo'
$ws.Range('B746').Value = 'San Salvador El Verde'
$ws.Range('B748').Value = 'Tecali De Herrera'
$ws.Range('B752').Value = 'Tepanco De López'
$ws.Range('B754').Value = 'Tetela De Ocampo'
$ws.Range('B758').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B768').Value = 'Amealco De Bonfil'
$ws.Range('B769').Value = 'Cadereyta De Montes'
$ws.Range('B774').Value = 'Jalpan De Serra'
$ws.Range('B775').Value = 'Landa De Matamoros'
$ws.Range('B777').Value = 'Pinal De Amoles'
$ws.Range('B797').Value = 'Santa María Del Río'
$ws.Range('B801').Value = 'Villa De Ramos'
$ws.Range('B859').Value = 'Nacozari De García'
$ws.Range('B870').Value = 'San Pedro De La Cueva'
$ws.Range('B912').Value = 'Ziltlaltépec De Trinidad Sánchez Santos'
$ws.Range('B920').Value = 'Amatlán De Los Reyes'
$ws.Range('B930').Value = 'Cosamaloapan De Carpio'
$ws.Range('B937').Value = 'Hueyapan De Ocampo'
$ws.Range('B938').Value = 'Ignacio De La Llave'
$ws.Range('B940').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B947').Value = 'Lerdo De Tejada'
$ws.Range('B950').Value = 'Martínez De La Torre'
$ws.Range('B951').Value = 'Medellín De Bravo'
$ws.Range('B954').Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Range('B964').Value = 'Paso De Ovejas'
$ws.Range('B966').Value = 'Poza Rica De Hidalgo'
$ws.Range('B972').Value = 'Soledad De Doblado'
$ws.Range('B984').Value = 'Vega De Alatorre'
$ws.Range('B1000').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B1003').Value = 'El Plateado De Joaquín Amaro'
$ws.Range('B1010').Value = 'Jiménez Del Teul'
$ws.Range('B1019').Value = 'Moyahua De Estrada'
$ws.Range('B1020').Value = 'Nochistlán De Mejía'
$ws.Range('B1021').Value = 'Noria De Ángeles'
$ws.Range('B1030').Value = 'Teúl De González Ortega'
$ws.Range('B1031').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B1034').Value = 'Villa De Cos'
$ws.Range('A258').Value = 'Guanajuato'
$ws.Range('B616').Value = 'Montemorelos'

# Remove the trailing footnote rows (source/sample-size/attribution text) that followed the data
$ws.Rows("1041:1045").Delete()
